# B6-PowerPoint.pptx edit
#
# 1. The three data tables in the deck (slides 14, 15 and 16) get their
#    table style switched from the plain "Table_0" style
#    ({873A5A61-9964-4137-8DD6-1F4297E06702}) to the built-in style
#    {9708B5FA-0537-48FF-AD6F-CA0408741185}.
#
# 2. The presentation's theme colours are swapped from the custom
#    "Integral / Red Violet" palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style every table on the slides -------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle("{9708B5FA-0537-48FF-AD6F-CA0408741185}")
        }
    }
}

# --- 2. Swap the theme colour scheme back to the default Office palette --
$cs = $p.SlideMaster.Theme.ThemeColorScheme
$cs.Item(1).RGB  = 0x000000   # dk1
$cs.Item(2).RGB  = 0xFFFFFF   # lt1
$cs.Item(3).RGB  = 0x6A5444   # dk2      (44546A)
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      (E7E6E6)
$cs.Item(5).RGB  = 0xD59B5B   # accent1  (5B9BD5)
$cs.Item(6).RGB  = 0x317DED   # accent2  (ED7D31)
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  (A5A5A5)
$cs.Item(8).RGB  = 0x00C0FF   # accent4  (FFC000)
$cs.Item(9).RGB  = 0xC47244   # accent5  (4472C4)
$cs.Item(10).RGB = 0x47AD70   # accent6  (70AD47)
$cs.Item(11).RGB = 0xC16305   # hlink    (0563C1)
$cs.Item(12).RGB = 0x724F95   # folHlink (954F72)
